# One-click update from Van Paper 01:03 PM on 2025-09-12
#
# 1) Row 36 (BEIRUT LEBANESE STREET FOOD) gets a "Last Invoice Date" of 09/12/2025.
# 2) A new customer row "Shakopee Brewhall" (customer # 0008296, salesperson
#    Dack Suzanne, prospect code 015) is inserted as row 38, pushing the
#    previously-last four rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the missing Last Invoice Date for row 36 -------------------
# Copy the formatting (date number format / style) from an existing
# date cell in the same column so the new value picks up the same style,
# then overwrite with the actual date serial for 2025-09-12.
$ws.Range("D39").Copy($ws.Range("D36"))
$ws.Range("D36").Value = 45912

# --- 2. Insert the new "Shakopee Brewhall" row at position 38 --------------
$ws.Rows.Item(38).Insert()

# Restore the row height that Excel normally carries for data rows in this
# sheet (the inserted row otherwise comes back without explicit height).
$ws.Rows.Item(38).RowHeight = 13.05

# Column F is an always-blank trailing column; copy it from the row above so
# the new row keeps the same (unstyled) blank cell in F.
$ws.Range("F37").Copy($ws.Range("F38"))

$ws.Range("A38").Value = "Shakopee Brewhall"
$ws.Range("B38").Value = "Dack, Suzanne"
$ws.Range("C38").Value = "015"
$ws.Range("E38").Value = "0008296"
